# Added dynamic calculation logic: new "should we buy the player" rows in the
# Buy/Sell Decision section (Tab 3), plus formula1 updates for the existing
# "Player club score" / "Player goal score" rows that feed the new decision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 19: "Is the player nice" -------------------------------------
$ws.Range("A19").Value = "Tab 3"
$ws.Range("B19").Value = "Buy/Sell Decision"
$ws.Range("C19").Value = "Is the player nice"
$ws.Range("D19").Value = "Niceness scale of player"
$ws.Range("E19").Value = "y"
$ws.Range("F19").Value = "dropdown"
$ws.Range("G19").Value = "Yes, No"
$ws.Range("I19").Value = "Yes"
$ws.Range("K19").Value = "n"
$ws.Range("P19").Value = "n"

# --- New row 20: "Has the player settled in life" -------------------------
$ws.Range("A20").Value = "Tab 3"
$ws.Range("B20").Value = "Buy/Sell Decision"
$ws.Range("C20").Value = "Has the player settled in life "
$ws.Range("D20").Value = "Talk about player's personal life (marriage, kids)"
$ws.Range("E20").Value = "y"
$ws.Range("F20").Value = "dropdown"
$ws.Range("G20").Value = "Yes, No"
$ws.Range("I20").Value = "Yes"
$ws.Range("K20").Value = "n"
$ws.Range("P20").Value = "n"

# --- New row 21: "Is the player around his prime" -------------------------
$ws.Range("A21").Value = "Tab 3"
$ws.Range("B21").Value = "Buy/Sell Decision"
$ws.Range("C21").Value = "Is the player around his prime"
$ws.Range("D21").Value = "Near his prime in career in club and country football"
$ws.Range("E21").Value = "y"
$ws.Range("F21").Value = "dropdown"
$ws.Range("G21").Value = "Yes, No"
$ws.Range("I21").Value = "Yes"
$ws.Range("K21").Value = "n"
$ws.Range("P21").Value = "n"

# --- New row 22: "should we buy the player" (computed, read-only) ---------
$ws.Range("A22").Value = "Tab 3"
$ws.Range("B22").Value = "Buy/Sell Decision"
$ws.Range("C22").Value = "should we buy the player"
$ws.Range("E22").Value = "e"
$ws.Range("F22").Value = "read-only"
$ws.Range("H22").Value = "yesand(Does the player have good fan following, Is the player nice, Is the player around his prime)"

# --- Update existing formula1 / field-name text on rows 15 & 16 -----------
# Row 15 becomes the "Player goal score" row, with its calc formula text.
$ws.Range("H15").Value = "Goals scored for country + Goals scored for club + Goals scored in youth career"
$ws.Range("C15").Value = "Player goal score"

# Row 16 ("Player club score") gets its calc formula text.
$ws.Range("H16").Value = "Goals scored for club/No. of clubs played for"

# --- Resize the Table1 ListObject (and AutoFilter) to cover new rows ------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:S22"))

# --- Update view state: scroll + active selection --------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
